# Auto-generated Excel COM-interop script to apply Golem_Profits numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 571.1667
$ws.Range("I19").Value = 649.6
$ws.Range("K19").Value = 649.6
$ws.Range("M19").Value = -474.6
$ws.Range("H33").Value = 852.8
$ws.Range("I33").Value = 879
$ws.Range("J33").Value = 748
$ws.Range("K33").Value = 879
$ws.Range("L33").Value = 748
$ws.Range("M33").Value = -650
$ws.Range("N33").Value = -1206
$ws.Range("H40").Value = 2309.8
$ws.Range("J40").Value = 2356.8572
$ws.Range("L40").Value = 2356.8572
$ws.Range("N40").Value = -2706.8572
$ws.Range("H53").Value = 138.8
$ws.Range("I53").Value = 40
$ws.Range("K53").Value = 40
$ws.Range("M53").Value = 597
$ws.Range("H55").Value = 2898.5715
$ws.Range("I55").Value = 2965
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 2965
$ws.Range("L55").Value = 2500
$ws.Range("M55").Value = -2751
$ws.Range("N55").Value = -2928
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("H70").Value = 2175.7144
$ws.Range("I70").Value = 1750
$ws.Range("J70").Value = 2346
$ws.Range("K70").Value = 5250
$ws.Range("L70").Value = 7038
$ws.Range("M70").Value = -4980
$ws.Range("N70").Value = -7578
$ws.Range("H73").Value = 2175.7144
$ws.Range("I73").Value = 1750
$ws.Range("J73").Value = 2346
$ws.Range("K73").Value = 5250
$ws.Range("L73").Value = 7038
$ws.Range("M73").Value = -4314
$ws.Range("N73").Value = -8910
$ws.Range("H98").Value = 692.375
$ws.Range("I98").Value = 876.8
$ws.Range("J98").Value = 385
$ws.Range("K98").Value = 876.8
$ws.Range("L98").Value = 385
$ws.Range("M98").Value = 621.2
$ws.Range("N98").Value = -3381
$ws.Range("H107").Value = 43237.094
$ws.Range("I107").Value = 53309.766
$ws.Range("J107").Value = 428.25
$ws.Range("K107").Value = 53309.766
$ws.Range("L107").Value = 428.25
$ws.Range("M107").Value = -51389.766
$ws.Range("N107").Value = -4268.25
$ws.Range("H122").Value = 692.375
$ws.Range("I122").Value = 876.8
$ws.Range("J122").Value = 385
$ws.Range("K122").Value = 2630.4
$ws.Range("L122").Value = 1155
$ws.Range("M122").Value = -180.3999999999996
$ws.Range("N122").Value = -6055
$ws.Range("H138").Value = 2711.2778
$ws.Range("J138").Value = 3170.1538
$ws.Range("L138").Value = 9510.4614
$ws.Range("N138").Value = -19790.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 676.5
$ws.Range("I2").Value = 711.8
$ws.Range("K2").Value = 711.8
$ws.Range("M2").Value = -598.8
$ws.Range("H116").Value = 676.5
$ws.Range("I116").Value = 711.8
$ws.Range("K116").Value = 711.8
$ws.Range("M116").Value = 1582.2
$ws.Range("H118").Value = 25500
$ws.Range("J118").Value = 25500
$ws.Range("L118").Value = 25500
$ws.Range("N118").Value = -28814

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 676.5
$ws.Range("I3").Value = 711.8
$ws.Range("K3").Value = 711.8
$ws.Range("M3").Value = -597.8
$ws.Range("H7").Value = 2250.5
$ws.Range("I7").Value = 1001
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 1001
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -888
$ws.Range("N7").Value = -3726
$ws.Range("H20").Value = 1100.7142
$ws.Range("I20").Value = 923.25
$ws.Range("J20").Value = 1337.3334
$ws.Range("K20").Value = 923.25
$ws.Range("L20").Value = 1337.3334
$ws.Range("M20").Value = -676.25
$ws.Range("N20").Value = -1831.3334
$ws.Range("H22").Value = 810.3333
$ws.Range("I22").Value = 756.4286
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 756.4286
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -583.4286
$ws.Range("N22").Value = -1345
$ws.Range("H105").Value = 893.9091
$ws.Range("I105").Value = 826.2222
$ws.Range("J105").Value = 1198.5
$ws.Range("K105").Value = 826.2222
$ws.Range("L105").Value = 1198.5
$ws.Range("M105").Value = 920.7778
$ws.Range("N105").Value = -4692.5
$ws.Range("H107").Value = 134816.33
$ws.Range("I107").Value = 200500
$ws.Range("J107").Value = 3449
$ws.Range("K107").Value = 200500
$ws.Range("L107").Value = 3449
$ws.Range("M107").Value = -198580
$ws.Range("N107").Value = -7289
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 107.07692
$ws.Range("I7").Value = 235.8
$ws.Range("K7").Value = 235.8
$ws.Range("M7").Value = -122.8
$ws.Range("H64").Value = 43750
$ws.Range("I64").Value = 15000
$ws.Range("J64").Value = 53333.332
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = 53333.332
$ws.Range("M64").Value = -14752
$ws.Range("N64").Value = -53829.332
$ws.Range("H67").Value = 43750
$ws.Range("I67").Value = 15000
$ws.Range("J67").Value = 53333.332
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 53333.332
$ws.Range("M67").Value = -14142
$ws.Range("N67").Value = -55049.332
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = $null
$ws.Range("H96").Value = 21047.2
$ws.Range("J96").Value = 21047.2
$ws.Range("L96").Value = 21047.2
$ws.Range("N96").Value = -26539.2
$ws.Range("H141").Value = 845999
$ws.Range("J141").Value = 1383332
$ws.Range("L141").Value = 1383332
$ws.Range("N141").Value = -1393692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4628.467
$ws.Range("I4").Value = 6045.952
$ws.Range("J4").Value = 1321
$ws.Range("K4").Value = 18137.856
$ws.Range("L4").Value = 3963
$ws.Range("M4").Value = -18025.856
$ws.Range("N4").Value = -4187
$ws.Range("H113").Value = 524.375
$ws.Range("I113").Value = 633.6667
$ws.Range("J113").Value = 196.5
$ws.Range("K113").Value = 1901.0001
$ws.Range("L113").Value = 589.5
$ws.Range("M113").Value = 268.9999
$ws.Range("N113").Value = -4929.5
$ws.Range("H131").Value = 2867
$ws.Range("I131").Value = 830
$ws.Range("J131").Value = 4904
$ws.Range("K131").Value = 2490
$ws.Range("L131").Value = 14712
$ws.Range("M131").Value = 2550
$ws.Range("N131").Value = -24792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 100000
$ws.Range("I10").Value = 100000
$ws.Range("K10").Value = 100000
$ws.Range("M10").Value = -99831
$ws.Range("H107").Value = 19608504
$ws.Range("I107").Value = 192.9
$ws.Range("J107").Value = 47620376
$ws.Range("K107").Value = 192.9
$ws.Range("L107").Value = 47620376
$ws.Range("M107").Value = 1727.1
$ws.Range("N107").Value = -47624216

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3005
$ws.Range("I20").Value = 3005
$ws.Range("K20").Value = 3005
$ws.Range("M20").Value = -2779
$ws.Range("H40").Value = 1390
$ws.Range("I40").Value = 1390
$ws.Range("K40").Value = 1390
$ws.Range("M40").Value = -1254
$ws.Range("H42").Value = 33999.8
$ws.Range("J42").Value = 33999.8
$ws.Range("L42").Value = 33999.8
$ws.Range("N42").Value = -35125.8
$ws.Range("H49").Value = 33999.8
$ws.Range("J49").Value = 33999.8
$ws.Range("L49").Value = 33999.8
$ws.Range("N49").Value = -34293.8
$ws.Range("H55").Value = 527.82355
$ws.Range("I55").Value = 385.66666
$ws.Range("J55").Value = 869
$ws.Range("K55").Value = 385.66666
$ws.Range("L55").Value = 869
$ws.Range("M55").Value = -212.66666
$ws.Range("N55").Value = -1215
$ws.Range("H132").Value = 2880
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2579.25
$ws.Range("I2").Value = 2579.25
$ws.Range("K2").Value = 2579.25
$ws.Range("M2").Value = -2467.25
$ws.Range("H62").Value = 4651.3335
$ws.Range("I62").Value = 4651.3335
$ws.Range("K62").Value = 4651.3335
$ws.Range("M62").Value = -4027.3335
$ws.Range("H65").Value = 4651.3335
$ws.Range("I65").Value = 4651.3335
$ws.Range("K65").Value = 23256.6675
$ws.Range("M65").Value = -20136.6675
$ws.Range("H126").Value = 3007
$ws.Range("I126").Value = 2550.7334
$ws.Range("J126").Value = 3862.5
$ws.Range("K126").Value = 7652.2002
$ws.Range("L126").Value = 11587.5
$ws.Range("M126").Value = -5182.2002
$ws.Range("N126").Value = -16527.5
$ws.Range("H132").Value = 1645
$ws.Range("I132").Value = 1225
$ws.Range("K132").Value = 3675
$ws.Range("M132").Value = -1145
